# Plantilla_SCORM: add an "Correo electrónico" column after Id_curso,
# and add "Elemento" / "Otra columna" columns at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before the old "Nombre" column (old B), pushing
#    Nombre..Otra cosa from B:G to C:H.
$ws.Columns("B:B").Insert()

# 2) Insert two new columns after the (now shifted) "Otra cosa" column (H),
#    to host "Elemento" and "Otra columna".
$ws.Columns("I:J").Insert()

# --- suspend_data column (now F) becomes the same value on every row ---
# (written first, matching the original authoring order)
$suspendData = "00119.4Ā-9Ą1720;āċ071812ĆĖ47Đěē8č39Ġ902đĞģģ03ħĔĩīĤĂĤĞĔ3ĕ8ĶĥăĒĔĵķĹĲļĿ9ĔĆ04ĻĴ8Ęŋ4"
$ws.Range("F2").Value = $suspendData
$ws.Range("F3").Value = $suspendData
$ws.Range("F4").Value = $suspendData
$ws.Range("F5").Value = $suspendData

# --- New "Elemento" column ---
$ws.Range("I1").Value = "Elemento"
$ws.Range("I2").Value = "Sí"
$ws.Range("I4").Value = "No"
$ws.Range("I5").Value = "Tal vez"

# --- New "Correo electrónico" column ---
$ws.Range("B1").Value = "Correo electrónico"
$ws.Range("B2").Value = "mario"
$ws.Range("B3").Value = "omar"
$ws.Range("B4").Value = "fer"
$ws.Range("B5").Value = "isma"

# --- New "Otra columna" column ---
$ws.Range("J1").Value = "Otra columna"
$ws.Range("J2").Value = 8
$ws.Range("J3").Value = 9
$ws.Range("J4").Value = 6
$ws.Range("J5").Value = 5

# --- Column widths (best-fit, matches final layout) ---
$ws.Columns.Item(1).ColumnWidth = 7.666666666666667
$ws.Columns.Item(2).ColumnWidth = 7.666666666666667
$ws.Columns.Item(3).ColumnWidth = 8.666666666666666
$ws.Columns.Item(4).ColumnWidth = 10.333333333333334
$ws.Columns.Item(5).ColumnWidth = 10.333333333333334
$ws.Columns.Item(6).ColumnWidth = 81.5
$ws.Columns.Item(7).ColumnWidth = 6.666666666666667
$ws.Columns.Item(8).ColumnWidth = 8.333333333333334

# --- Selection matches the post-edit workbook ---
[void]$ws.Range("J6").Select()

Write-Output "edit applied"
